$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024-12-19 00:22:13", -0.1214307487694896, -0.001914794722293995, 0.009300598274721067),
    @("2024-12-19 00:22:14", -0.122357947009732, -0.001832193056181995, 0.008967335235196622),
    @("2024-12-19 00:22:15", -0.1215771484916331, -0.001779453512009995, 0.008653635354543891),
    @("2024-12-19 00:22:16", -0.1215381085657282, -0.001856993801521996, 0.009027820570210583),
    @("2024-12-19 00:22:17", -0.1219773077321588, -0.001865496914209996, 0.009101931647119412),
    @("2024-12-19 00:22:18", -0.1236690378547064, -0.001842923174573997, 0.009116501433588293),
    @("2024-12-19 00:22:19", -0.1242448767618043, -0.001668862433299996, 0.008293904294310518),
    @("2024-12-19 00:22:20", -0.1242123434902169, -0.001706215392607996, 0.008477320496596785)
)

$startRow = 127
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
